# Apply the "Список заданий" task-list rewrite described by the commit.
#
# Strategy: for paragraphs whose whole run content changes we rebuild the
# paragraph's run content from scratch via Range.InsertXML (a minimal
# WordprocessingML package fragment) - this lets us control run
# boundaries precisely (needed for the w:proofErr-wrapped misspelled
# gaming terms) instead of fighting Range.Text's single-run semantics.
#
# NOTE: this interpreter parses `$comObject (expr)` (a COM-object variable
# immediately followed by a parenthesised expression) as an *invocation*
# of that variable, which silently corrupts sibling arguments in the same
# command call. To stay safe, every helper-function result that will be
# passed alongside a paragraph/range object is first stashed in a plain
# temp variable, and only bare variables are passed as call arguments.

$d = $word.ActiveDocument

function New-RunXml($text, [bool]$preserve = $false) {
    $space = ""
    if ($preserve) { $space = ' xml:space="preserve"' }
    return '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t' + $space + '>' + $text + '</w:t></w:r>'
}

function New-ProofRunXml($text) {
    $run = New-RunXml $text
    return '<w:proofErr w:type="spellStart"/>' + $run + '<w:proofErr w:type="spellEnd"/>'
}

function Wrap-ParaPkg($innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replace a whole paragraph's run content (everything except the trailing
# paragraph mark) with $runsXml (one or more <w:r>/<w:proofErr> elements),
# optionally prefixed with $prefixXml (e.g. a bookmark pair).
function Set-ParaRuns($para, $runsXml, $prefixXml) {
    $r = $para.Range
    $start = $r.Start
    $end = $r.End
    $target = $d.Range($start, $end - 1)
    $innerPara = '<w:p>' + $prefixXml + $runsXml + '</w:p>'
    $pkg = Wrap-ParaPkg $innerPara
    $target.InsertXML($pkg)
}

# --- 1. Drop the _GoBack bookmark from the title paragraph --------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- 2. Drop the three sub-task paragraphs (old "А) ..." / "Б) ..." /
#        "атаковать ..." block) that followed the old item 6 -------------
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Delete()
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Delete()
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Delete()

# --- 3. Rewrite the remaining list items in place ------------------------

# Item 3: old multi-run "Запилить ИИ 2.0. + ..." -> bookmark + new text
$bookmarkXml = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$runs3 = New-RunXml "Сделать так, чтобы заход и выход из меню паузы никак не выводило из паузы, связанной со стартовым диалогом."
$item3 = $d.Paragraphs.Item(3)
Set-ParaRuns $item3 $runs3 $bookmarkXml

# Item 4: -> "Создать такой тип хитбоксов, который бьёт и героя, и врагов."
$run4a = New-RunXml "Создать такой тип " $true
$run4b = New-ProofRunXml "хитбоксов"
$run4c = New-RunXml ", который бьёт и героя, и врагов."
$runs4 = $run4a + $run4b + $run4c
$noPrefix = ""
$item4 = $d.Paragraphs.Item(4)
Set-ParaRuns $item4 $runs4 $noPrefix

# Item 5: -> "Сделать игру весёлой"
$runs5 = New-RunXml "Сделать игру весёлой"
$item5 = $d.Paragraphs.Item(5)
Set-ParaRuns $item5 $runs5 $noPrefix

# Item 6: -> "Поработать над корректностью размеров окошек ГУИ."
$runs6 = New-RunXml "Поработать над корректностью размеров окошек ГУИ."
$item6 = $d.Paragraphs.Item(6)
Set-ParaRuns $item6 $runs6 $noPrefix

# Item 7 (was "Добавить новых монстров") ->
#   "Исправить джунглиевого паука, гигантскую летучую мышь, босса шахтёров-призраков"
$run7a = New-RunXml "Исправить " $true
$run7b = New-ProofRunXml "джунглиевого"
$run7c = New-RunXml " паука, гигантскую летучую мышь, босса шахтёров-призраков" $true
$runs7 = $run7a + $run7b + $run7c
$item7 = $d.Paragraphs.Item(7)
Set-ParaRuns $item7 $runs7 $noPrefix

# Item 8 (was "Сделать игру весёлой") ->
#   "Добавить события, которые влияют на всю игру, а не только на в пределах уровня"
$runs8 = New-RunXml "Добавить события, которые влияют на всю игру, а не только на в пределах уровня"
$item8 = $d.Paragraphs.Item(8)
Set-ParaRuns $item8 $runs8 $noPrefix

# Item 9 (was "Поработать над корректностью размеров окошек ГУИ.") ->
#   "Добавить кнопку выхода"
$runs9 = New-RunXml "Добавить кнопку выхода"
$item9 = $d.Paragraphs.Item(9)
Set-ParaRuns $item9 $runs9 $noPrefix

# --- 4. Append four brand-new list items ---------------------------------
$pPrXml = '<w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'

$runsN1 = New-RunXml "Добавить дрожание экрана при ударах и получении урона."
$newPara1 = '<w:p>' + $pPrXml + $runsN1 + '</w:p>'

$runsN2 = New-RunXml "Добавить музыку и звуки."
$newPara2 = '<w:p>' + $pPrXml + $runsN2 + '</w:p>'

$runsN3 = New-RunXml "Убрать дёргания при перемещении союзника."
$newPara3 = '<w:p>' + $pPrXml + $runsN3 + '</w:p>'

$runN4a = New-RunXml "Убрать " $true
$runN4b = New-ProofRunXml "читы"
$runsN4 = $runN4a + $runN4b
$newPara4 = '<w:p>' + $pPrXml + $runsN4 + '</w:p>'

$lastIndex = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($lastIndex)
$lastEnd = $last.Range.End
$insertPoint = $d.Range($lastEnd, $lastEnd)
$newBody = $newPara1 + $newPara2 + $newPara3 + $newPara4
$pkg = Wrap-ParaPkg $newBody
$insertPoint.InsertXML($pkg)

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
